$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The updated data only retains four target clusters per sending
# cluster (ECs, FAPs, Inflammatory-Mac, MuSCs); the 'Neutrophils' and
# 'Resolving-Mac' target-cluster rows are removed for both the ECs-
# sender block (originally rows 6-7) and the FAPs-sender block
# (originally rows 12-13, i.e. rows 10-11 after the first deletion).
$ws.Rows.Item(6).Delete()
$ws.Rows.Item(6).Delete()
$ws.Rows.Item(10).Delete()
$ws.Rows.Item(10).Delete()

# Refresh the TPM-derived numeric columns (E:T) for the 8 remaining
# data rows with the newly computed values.
$ws.Range("E2").Value = 1
$ws.Range("F2").Value = 0.5
$ws.Range("G2").Value = 0.0292345
$ws.Range("H2").Value = 0.058469
$ws.Range("I2").Value = 0.4428765120700495
$ws.Range("J2").Value = 0.346386487911515
$ws.Range("K2").Value = 2
$ws.Range("L2").Value = 1
$ws.Range("M2").Value = 0.1557005
$ws.Range("N2").Value = 0.311401
$ws.Range("O2").Value = 0.02102398211576467
$ws.Range("P2").Value = 0.01500040222529337
$ws.Range("Q2").Value = 0.00455182626725
$ws.Range("R2").Value = 0.018207305069
$ws.Range("S2").Value = 0.009311027869252954
$ws.Range("T2").Value = 0.005195936644079444

$ws.Range("E3").Value = 1
$ws.Range("F3").Value = 0.5
$ws.Range("G3").Value = 0.0292345
$ws.Range("H3").Value = 0.058469
$ws.Range("I3").Value = 0.4428765120700495
$ws.Range("J3").Value = 0.346386487911515
$ws.Range("K3").Value = 3
$ws.Range("L3").Value = 1
$ws.Range("M3").Value = 5.923689
$ws.Range("N3").Value = 17.771067
$ws.Range("O3").Value = 0.7998659708565604
$ws.Range("P3").Value = 0.8560446272575798
$ws.Range("Q3").Value = 0.1731760860705
$ws.Range("R3").Value = 1.039056516423
$ws.Range("S3").Value = 0.3542418512964773
$ws.Range("T3").Value = 0.2965222919312751

$ws.Range("E4").Value = 1
$ws.Range("F4").Value = 0.5
$ws.Range("G4").Value = 0.0292345
$ws.Range("H4").Value = 0.058469
$ws.Range("I4").Value = 0.4428765120700495
$ws.Range("J4").Value = 0.346386487911515
$ws.Range("K4").Value = 1
$ws.Range("L4").Value = 0.3333333333333333
$ws.Range("M4").Value = 0.024117
$ws.Range("N4").Value = 0.072351
$ws.Range("O4").Value = 0.003256478795417461
$ws.Range("P4").Value = 0.003485197868350457
$ws.Range("Q4").Value = 0.0007050484365
$ws.Range("R4").Value = 0.004230290619
$ws.Range("S4").Value = 0.001442217970544561
$ws.Range("T4").Value = 0.001207225449294613

$ws.Range("E5").Value = 1
$ws.Range("F5").Value = 0.5
$ws.Range("G5").Value = 0.0292345
$ws.Range("H5").Value = 0.058469
$ws.Range("I5").Value = 0.4428765120700495
$ws.Range("J5").Value = 0.346386487911515
$ws.Range("K5").Value = 2
$ws.Range("L5").Value = 1
$ws.Range("M5").Value = 1.3023455
$ws.Range("N5").Value = 2.604691
$ws.Range("O5").Value = 0.1758535682322574
$ws.Range("P5").Value = 0.1254697726487764
$ws.Range("Q5").Value = 0.03807341951975
$ws.Range("R5").Value = 0.152293678079
$ws.Range("S5").Value = 0.0778814149337746
$ws.Range("T5").Value = 0.04346103388686591

$ws.Range("E6").Value = 1
$ws.Range("F6").Value = 0.3333333333333333
$ws.Range("G6").Value = 0.036776
$ws.Range("H6").Value = 0.110328
$ws.Range("I6").Value = 0.5571234879299505
$ws.Range("J6").Value = 0.6536135120884849
$ws.Range("K6").Value = 2
$ws.Range("L6").Value = 1
$ws.Range("M6").Value = 0.1557005
$ws.Range("N6").Value = 0.311401
$ws.Range("O6").Value = 0.02102398211576467
$ws.Range("P6").Value = 0.01500040222529337
$ws.Range("Q6").Value = 0.005726041587999999
$ws.Range("R6").Value = 0.034356249528
$ws.Range("S6").Value = 0.01171295424651171
$ws.Range("T6").Value = 0.009804465581213924

$ws.Range("E7").Value = 1
$ws.Range("F7").Value = 0.3333333333333333
$ws.Range("G7").Value = 0.036776
$ws.Range("H7").Value = 0.110328
$ws.Range("I7").Value = 0.5571234879299505
$ws.Range("J7").Value = 0.6536135120884849
$ws.Range("K7").Value = 3
$ws.Range("L7").Value = 1
$ws.Range("M7").Value = 5.923689
$ws.Range("N7").Value = 17.771067
$ws.Range("O7").Value = 0.7998659708565604
$ws.Range("P7").Value = 0.8560446272575798
$ws.Range("Q7").Value = 0.217849586664
$ws.Range("R7").Value = 1.960646279976
$ws.Range("S7").Value = 0.4456241195600831
$ws.Range("T7").Value = 0.5595223353263047

$ws.Range("E8").Value = 1
$ws.Range("F8").Value = 0.3333333333333333
$ws.Range("G8").Value = 0.036776
$ws.Range("H8").Value = 0.110328
$ws.Range("I8").Value = 0.5571234879299505
$ws.Range("J8").Value = 0.6536135120884849
$ws.Range("K8").Value = 1
$ws.Range("L8").Value = 0.3333333333333333
$ws.Range("M8").Value = 0.024117
$ws.Range("N8").Value = 0.072351
$ws.Range("O8").Value = 0.003256478795417461
$ws.Range("P8").Value = 0.003485197868350457
$ws.Range("Q8").Value = 0.0008869267919999999
$ws.Range("R8").Value = 0.007982341128
$ws.Range("S8").Value = 0.0018142608248729
$ws.Range("T8").Value = 0.002277972419055843

$ws.Range("E9").Value = 1
$ws.Range("F9").Value = 0.3333333333333333
$ws.Range("G9").Value = 0.036776
$ws.Range("H9").Value = 0.110328
$ws.Range("I9").Value = 0.5571234879299505
$ws.Range("J9").Value = 0.6536135120884849
$ws.Range("K9").Value = 2
$ws.Range("L9").Value = 1
$ws.Range("M9").Value = 1.3023455
$ws.Range("N9").Value = 2.604691
$ws.Range("O9").Value = 0.1758535682322574
$ws.Range("P9").Value = 0.1254697726487764
$ws.Range("Q9").Value = 0.047895058108
$ws.Range("R9").Value = 0.287370348648
$ws.Range("S9").Value = 0.09797215329848276
$ws.Range("T9").Value = 0.08200873876191045
